# Edit script generated to apply the commit diff to horarios-141-completo.xlsx
# Adds new scraped schedule rows (31/12/2025 16:58:55 / 16:59:01 scrape batches)
# to all three worksheets and refreshes the 'last updated' / 'total rows' header cells.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Update header/summary cells on each sheet ---
$ws1.Range("A2").Value = "Última actualización: 31/12/2025 16:59:06"
$ws1.Range("A3").Value = "Total filas: 1143"

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 16:59:06"
$ws2.Range("A3").Value = "Total filas: 76"

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 16:59:06"
$ws3.Range("A3").Value = "Total filas: 138"

# --- Sheet1 (LP1912): append 20 new rows 1125-1144 ---
$data = New-Object 'object[,]' 20,7
$data[0,0] = ""
$data[0,1] = '16:58:55'
$data[0,2] = '17:01'
$data[0,3] = '14_ABASTO'
$data[0,4] = 3
$data[0,5] = 'LP1912'
$data[0,6] = '31/12/2025'
$data[1,0] = ""
$data[1,1] = '16:58:55'
$data[1,2] = '17:03'
$data[1,3] = '23_HERNANDEZ'
$data[1,4] = 5
$data[1,5] = 'LP1912'
$data[1,6] = '31/12/2025'
$data[2,0] = ""
$data[2,1] = '16:58:55'
$data[2,2] = '17:07'
$data[2,3] = '15_ABASTO'
$data[2,4] = 9
$data[2,5] = 'LP1912'
$data[2,6] = '31/12/2025'
$data[3,0] = ""
$data[3,1] = '16:58:55'
$data[3,2] = '17:13'
$data[3,3] = '10_OLMOS'
$data[3,4] = 15
$data[3,5] = 'LP1912'
$data[3,6] = '31/12/2025'
$data[4,0] = ""
$data[4,1] = '16:58:55'
$data[4,2] = '17:23'
$data[4,3] = '16_SANTA ANA'
$data[4,4] = 25
$data[4,5] = 'LP1912'
$data[4,6] = '31/12/2025'
$data[5,0] = ""
$data[5,1] = '16:58:55'
$data[5,2] = '17:24'
$data[5,3] = '11_ETCHEVERRY'
$data[5,4] = 26
$data[5,5] = 'LP1912'
$data[5,6] = '31/12/2025'
$data[6,0] = ""
$data[6,1] = '16:58:55'
$data[6,2] = '17:27'
$data[6,3] = '15_ABASTO'
$data[6,4] = 29
$data[6,5] = 'LP1912'
$data[6,6] = '31/12/2025'
$data[7,0] = ""
$data[7,1] = '16:58:55'
$data[7,2] = '17:33'
$data[7,3] = '23_HERNANDEZ'
$data[7,4] = 35
$data[7,5] = 'LP1912'
$data[7,6] = '31/12/2025'
$data[8,0] = ""
$data[8,1] = '16:58:55'
$data[8,2] = '17:34'
$data[8,3] = '10_OLMOS'
$data[8,4] = 36
$data[8,5] = 'LP1912'
$data[8,6] = '31/12/2025'
$data[9,0] = ""
$data[9,1] = '16:58:55'
$data[9,2] = '17:35'
$data[9,3] = '16_P MOR-SANTA ANA'
$data[9,4] = 37
$data[9,5] = 'LP1912'
$data[9,6] = '31/12/2025'
$data[10,0] = ""
$data[10,1] = '16:58:55'
$data[10,2] = '17:38'
$data[10,3] = '17X38_ROMERO'
$data[10,4] = 40
$data[10,5] = 'LP1912'
$data[10,6] = '31/12/2025'
$data[11,0] = ""
$data[11,1] = '16:58:55'
$data[11,2] = '17:47'
$data[11,3] = '16_SANTA ANA'
$data[11,4] = 49
$data[11,5] = 'LP1912'
$data[11,6] = '31/12/2025'
$data[12,0] = ""
$data[12,1] = '16:58:55'
$data[12,2] = '17:50'
$data[12,3] = '215_EL PELIGRO'
$data[12,4] = 52
$data[12,5] = 'LP1912'
$data[12,6] = '31/12/2025'
$data[13,0] = ""
$data[13,1] = '16:58:55'
$data[13,2] = '17:54'
$data[13,3] = '10_OLMOS'
$data[13,4] = 56
$data[13,5] = 'LP1912'
$data[13,6] = '31/12/2025'
$data[14,0] = ""
$data[14,1] = '16:58:55'
$data[14,2] = '17:59'
$data[14,3] = '16_SANTA ANA'
$data[14,4] = 61
$data[14,5] = 'LP1912'
$data[14,6] = '31/12/2025'
$data[15,0] = ""
$data[15,1] = '16:58:55'
$data[15,2] = '18:01'
$data[15,3] = '23_HERNANDEZ'
$data[15,4] = 63
$data[15,5] = 'LP1912'
$data[15,6] = '31/12/2025'
$data[16,0] = ""
$data[16,1] = '16:58:55'
$data[16,2] = '18:04'
$data[16,3] = '14_ABASTO'
$data[16,4] = 66
$data[16,5] = 'LP1912'
$data[16,6] = '31/12/2025'
$data[17,0] = ""
$data[17,1] = '16:58:55'
$data[17,2] = '18:24'
$data[17,3] = '11_ETCHEVERRY'
$data[17,4] = 86
$data[17,5] = 'LP1912'
$data[17,6] = '31/12/2025'
$data[18,0] = ""
$data[18,1] = '16:58:55'
$data[18,2] = '18:27'
$data[18,3] = '15_ABASTO'
$data[18,4] = 89
$data[18,5] = 'LP1912'
$data[18,6] = '31/12/2025'
$data[19,0] = ""
$data[19,1] = '16:58:55'
$data[19,2] = '18:33'
$data[19,3] = '14X44_ABASTO'
$data[19,4] = 95
$data[19,5] = 'LP1912'
$data[19,6] = '31/12/2025'
$rng = $ws1.Range("A1125:G1144")
$rng.Value = $data

# --- Sheet2 (LP1912-215): append 1 new row 77 ---
$data = New-Object 'object[,]' 1,7
$data[0,0] = ""
$data[0,1] = '31/12/2025'
$data[0,2] = '16:58:55'
$data[0,3] = '17:50'
$data[0,4] = '215_EL PELIGRO'
$data[0,5] = 52
$data[0,6] = 'LP1912'
$rng = $ws2.Range("A77:G77")
$rng.Value = $data

# --- Sheet3 (6203-6173): append 2 new rows 138-139 ---
$data = New-Object 'object[,]' 2,7
$data[0,0] = ""
$data[0,1] = '31/12/2025'
$data[0,2] = '16:59:01'
$data[0,3] = '17:01'
$data[0,4] = '215C_LA PLATA'
$data[0,5] = 2
$data[0,6] = 'L6203'
$data[1,0] = ""
$data[1,1] = '31/12/2025'
$data[1,2] = '16:59:01'
$data[1,3] = '18:22'
$data[1,4] = '215C_LA PLATA'
$data[1,5] = 83
$data[1,6] = 'L6203'
$rng = $ws3.Range("A138:G139")
$rng.Value = $data

